$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.423.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.73%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.641.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.56%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'212.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.529"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'23.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.34%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.49%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.86%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.873.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.612.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.573"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.17%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.40%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.37%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'27.391.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.81%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'229.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.47%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.77%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.04%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.97%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.12%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'147.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +1.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.18%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.54%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -4.30%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.56%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.08%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.58%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.414.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.64%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.71%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.886"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.42%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'TrustWalletToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.805"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'MXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'64.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.38%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.783.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'87.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.68%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.20%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Algorand"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.87%  "
$ws.Range("E51").Style = "Normal"
